$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.8
$ws.Range("H2").Value = 4
$ws.Range("I2").Value = 3.8
$ws.Range("J2").Value = 2.4
$ws.Range("K2").Value = 2.38
$ws.Range("L2").Value = 4.33
$ws.Range("N2").Value = 15
$ws.Range("O2").Value = 1.18
$ws.Range("P2").Value = 4.5
$ws.Range("Q2").Value = 1.65
$ws.Range("R2").Value = 2.2
$ws.Range("S2").Value = 1.3
$ws.Range("T2").Value = 3.4
$ws.Range("U2").Value = 1.57
$ws.Range("V2").Value = 2.25
$ws.Range("W2").Value = 9.5
$ws.Range("X2").Value = 10
$ws.Range("Z2").Value = 15
$ws.Range("AB2").Value = 21
$ws.Range("AC2").Value = 15
$ws.Range("AE2").Value = 13
$ws.Range("AF2").Value = 41
$ws.Range("AG2").Value = 151
$ws.Range("AH2").Value = 15
$ws.Range("AJ2").Value = 13
$ws.Range("AN2").Value = 4
$ws.Range("AS2").Value = 101
$ws.Range("AT2").Value = 3.4
$ws.Range("AU2").Value = 7.5
$ws.Range("AY2").Value = 21
$ws.Range("BA2").Value = 67
$ws.Range("BC2").Value = 151

# Row 4
$ws.Range("G4").Value = 1.7
$ws.Range("H4").Value = 3.35
$ws.Range("I4").Value = 4.5
$ws.Range("J4").Value = 2.18
$ws.Range("K4").Value = 2.18
$ws.Range("L4").Value = 4.75
$ws.Range("M4").Value = 1.03
$ws.Range("N4").Value = 9.35
$ws.Range("O4").Value = 1.26
$ws.Range("P4").Value = 3.58
$ws.Range("Q4").Value = 1.82
$ws.Range("R4").Value = 1.8
$ws.Range("S4").Value = 1.43
$ws.Range("T4").Value = 2.75
$ws.Range("U4").Value = 1.81
$ws.Range("V4").Value = 1.92
$ws.Range("W4").Value = 6
$ws.Range("X4").Value = 6.9
$ws.Range("Z4").Value = 11.25
$ws.Range("AA4").Value = 11
$ws.Range("AC4").Value = 9.5
$ws.Range("AD4").Value = 5.8
$ws.Range("AE4").Value = 12
$ws.Range("AF4").Value = 50
$ws.Range("AH4").Value = 10.25
$ws.Range("AI4").Value = 21
$ws.Range("AJ4").Value = 12
$ws.Range("AK4").Value = 60
$ws.Range("AL4").Value = 35
$ws.Range("AM4").Value = 37
$ws.Range("AN4").Value = 3.6
$ws.Range("AO4").Value = 8
$ws.Range("AP4").Value = 15.5
$ws.Range("AQ4").Value = 26
$ws.Range("AS4").Value = 175
$ws.Range("AT4").Value = 2.65
$ws.Range("AU4").Value = 7
$ws.Range("AV4").Value = 55
$ws.Range("AX4").Value = 6.5
$ws.Range("AY4").Value = 26
$ws.Range("AZ4").Value = 29
$ws.Range("BA4").Value = 150
$ws.Range("BB4").Value = 175
$ws.Range("BC4").Value = 350

# Row 5
$ws.Range("G5").Value = 2.05
$ws.Range("I5").Value = 3
$ws.Range("J5").Value = 2.63
$ws.Range("L5").Value = 3.5
$ws.Range("M5").Value = 1.02
$ws.Range("O5").Value = 1.13
$ws.Range("S5").Value = 1.29
$ws.Range("T5").Value = 3.5
$ws.Range("Y5").Value = 9
$ws.Range("AA5").Value = 15
$ws.Range("AH5").Value = 15
$ws.Range("AI5").Value = 19
$ws.Range("AJ5").Value = 12
$ws.Range("AK5").Value = 34
$ws.Range("AM5").Value = 26
$ws.Range("AO5").Value = 11
$ws.Range("AP5").Value = 17
$ws.Range("AQ5").Value = 34
$ws.Range("AR5").Value = 41
$ws.Range("AT5").Value = 3.5
$ws.Range("AX5").Value = 5.5
$ws.Range("AY5").Value = 17
$ws.Range("BA5").Value = 51

# Row 9
$ws.Range("M9").Value = 1.1
$ws.Range("N9").Value = 7
$ws.Range("Q9").Value = 2.6
$ws.Range("R9").Value = 1.48

# Row 11
$ws.Range("G11").Value = 4
$ws.Range("H11").Value = 3.15
$ws.Range("I11").Value = 1.88
$ws.Range("J11").Value = 4.6
$ws.Range("K11").Value = 2.02
$ws.Range("L11").Value = 2.5
$ws.Range("M11").Value = 1.09
$ws.Range("N11").Value = 6.1
$ws.Range("O11").Value = 1.42
$ws.Range("P11").Value = 2.67
$ws.Range("Q11").Value = 2.22
$ws.Range("R11").Value = 1.6
$ws.Range("S11").Value = 1.47
$ws.Range("T11").Value = 2.52
$ws.Range("U11").Value = 2
$ws.Range("W11").Value = 9.5
$ws.Range("X11").Value = 21
$ws.Range("Y11").Value = 14
$ws.Range("Z11").Value = 65
$ws.Range("AA11").Value = 45
$ws.Range("AB11").Value = 55
$ws.Range("AC11").Value = 6.1
$ws.Range("AD11").Value = 6.2
$ws.Range("AE11").Value = 17
$ws.Range("AH11").Value = 5.9
$ws.Range("AI11").Value = 8
$ws.Range("AK11").Value = 16
$ws.Range("AL11").Value = 17
$ws.Range("AM11").Value = 35
$ws.Range("AN11").Value = 5.8
$ws.Range("AO11").Value = 24
$ws.Range("AP11").Value = 32
$ws.Range("AQ11").Value = 150
$ws.Range("AR11").Value = 200
$ws.Range("AT11").Value = 2.52
$ws.Range("AU11").Value = 7.7
$ws.Range("AX11").Value = 3.65
$ws.Range("AY11").Value = 9.75
$ws.Range("AZ11").Value = 21
$ws.Range("BA11").Value = 37
$ws.Range("BB11").Value = 80
